$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: previous day's date changes from 16/6/2025 (Mon) to 15/6/2025 (Sun)
$ws.Range("A2").Value = "15/6/2025 (Sun)"

# Row 3 gets this run's 4D box results
$ws.Range("A3").Value = "16/6/2025 (Mon)"
$ws.Range("B3").Value = "2 9 3 2`n0 2 8 9`n8 5 2 5`n6 7 4 1"
$ws.Range("C3").Value = "✅ Direct: 12/3547 (0.34%)`n✅ iBet: 12/195 (6.15%)"

# C6 becomes a blank cell matching the B/C formatted placeholder cells
$ws.Cells.Item(6, 3).WrapText = $true

# A new blank placeholder row (32) is appended at the bottom of the log
$ws.Cells.Item(32, 2).WrapText = $true
